$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

For ($row = 2; $row -le 123; $row++) {
    $ws.Cells.Item($row, 3).Value = 45180
}
